$d = $word.ActiveDocument

function Replace-In($paraIndex, $old, $new) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Paragraph 1: Title ---
Replace-In 1 'ContosoLearn Competitor SWOT' 'Analyse SWOT de la concurrence ContosoLearn'
# --- Paragraph 2: Fabrikam Learning heading ---
Replace-In 2 'Fabrikam Learning:' 'Fabrikam Learning :'
# --- Paragraph 3: Fabrikam - Strengths ---
Replace-In 3 'Strengths:' 'Forces :'
Replace-In 3 ' Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed.' ' Fabrikam Learning fournit un ensemble complet d’outils d’analyse et de création de rapports. Il garantit la surveillance continue des activités d’enseignement et d’apprentissage, ainsi que l’identification des domaines problématiques devant être abordés.'
# --- Paragraph 4: Fabrikam - Weaknesses ---
Replace-In 4 'Weaknesses:' 'Faiblesses :'
Replace-In 4 ' While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature.' ' bien que Fabrikam Learning dispose de fonctionnalités de création de rapports robustes, celles-ci peuvent être trop complexes pour certains utilisateurs, en raison de leur caractère exhaustif.'
# --- Paragraph 5: Fabrikam - Opportunities ---
Replace-In 5 'Opportunities:' 'Opportunités :'
Replace-In 5 ' There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand.' ' il y a une demande croissante d’expériences d’apprentissage personnalisées et de recommandations basées sur des données. Fabrikam Learning peut tirer parti de ses outils d’analytique et de création de rapports robustes pour répondre à cette demande.'
# --- Paragraph 6: Fabrikam - Threats ---
Replace-In 6 'Threats:' 'Menaces :'
Replace-In 6 ' The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead.' ' le marché de l’eLearning est très concurrentiel, avec de nombreux acteurs offrant des fonctionnalités similaires. Fabrikam Learning doit innover en permanence pour rester en tête.'
# --- Paragraph 7: AdatumLearn heading ---
Replace-In 7 'AdatumLearn:' 'AdatumLearn :'
# --- Paragraph 8: AdatumLearn - Strengths ---
Replace-In 8 'Strengths:' 'Points forts :'
Replace-In 8 ' AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users.' ' AdatumLearn propose des cours sur les techniques d’analyse métier telles que le MOST et le SWOT. Cela montre leur engagement à fournir du contenu précieux à leurs utilisateurs.'
# --- Paragraph 9: AdatumLearn - Weaknesses ---
Replace-In 9 'Weaknesses:' 'Faiblesses :'
Replace-In 9 ' The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content.' ' les informations fournies dans leurs cours sont un regroupement d’informations générées par des tiers. Celles-ci peuvent ne pas être aussi efficaces que du contenu original.'
# --- Paragraph 10: AdatumLearn - Opportunities ---
Replace-In 10 'Opportunities:' 'Opportunités :'
Replace-In 10 ' AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics.' ' AdatumLearn peut créer du contenu plus original pour fournir une valeur unique à leurs utilisateurs. Il peut également développer ses offres de cours afin de couvrir des sujets supplémentaires.'
# --- Paragraph 11: AdatumLearn - Threats ---
Replace-In 11 'Threats:' 'Menaces :'
Replace-In 11 ' Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive."' ' comme Fabrikam Learning, AdatumLearn fait également face à une rude concurrence sur le marché de l’eLearning. Il doit améliorer continuellement ses offres pour rester compétitif. »'

Write-Output "Done"
